$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3368.2432
$ws.Range("I132").Value = 3147.7942
$ws.Range("J132").Value = 5866.6665
$ws.Range("K132").Value = 9443.382599999999
$ws.Range("L132").Value = 17599.9995
$ws.Range("M132").Value = -6913.382599999999
$ws.Range("N132").Value = -22659.9995

$ws.Range("H138").Value = 2285.6155
$ws.Range("I138").Value = 1591.0526
$ws.Range("J138").Value = 2945.45
$ws.Range("K138").Value = 4773.1578
$ws.Range("L138").Value = 8836.349999999999
$ws.Range("M138").Value = 366.8422
$ws.Range("N138").Value = -19116.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7402.029
$ws.Range("I32").Value = 7738.912
$ws.Range("J32").Value = 5801.8335
$ws.Range("K32").Value = 7738.912
$ws.Range("L32").Value = 5801.8335
$ws.Range("M32").Value = -7451.912
$ws.Range("N32").Value = -6375.8335

$ws.Range("H61").Value = 41670176
$ws.Range("I61").Value = 50003870
$ws.Range("J61").Value = 1700
$ws.Range("K61").Value = 50003870
$ws.Range("L61").Value = 1700
$ws.Range("M61").Value = -50003658
$ws.Range("N61").Value = -2124

$ws.Range("H132").Value = 8335876.5
$ws.Range("I132").Value = 10418641
$ws.Range("J132").Value = 4818.6665
$ws.Range("K132").Value = 31255923
$ws.Range("L132").Value = 14455.9995
$ws.Range("M132").Value = -31253393
$ws.Range("N132").Value = -19515.9995

$ws.Range("H136").Value = 41670176
$ws.Range("I136").Value = 50003870
$ws.Range("J136").Value = 1700
$ws.Range("K136").Value = 150011610
$ws.Range("L136").Value = 5100
$ws.Range("M136").Value = -150009060
$ws.Range("N136").Value = -10200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10105140
$ws.Range("I31").Value = 4363.968
$ws.Range("J31").Value = 166667170
$ws.Range("K31").Value = 4363.968
$ws.Range("L31").Value = 166667170
$ws.Range("M31").Value = -4068.968

$ws.Range("H34").Value = 10105140
$ws.Range("I34").Value = 4363.968
$ws.Range("J34").Value = 166667170
$ws.Range("K34").Value = 4363.968
$ws.Range("L34").Value = 166667170
$ws.Range("M34").Value = -4161.968

$ws.Range("H58").Value = 1908.3793
$ws.Range("I58").Value = 741.05554
$ws.Range("J58").Value = 3818.5454
$ws.Range("K58").Value = 741.05554
$ws.Range("L58").Value = 3818.5454
$ws.Range("M58").Value = -538.05554
$ws.Range("N58").Value = -4224.5454

$ws.Range("H132").Value = 20006398
$ws.Range("I132").Value = 23815902
$ws.Range("J132").Value = 6503
$ws.Range("K132").Value = 71447706
$ws.Range("L132").Value = 19509
$ws.Range("M132").Value = -71445176
$ws.Range("N132").Value = -24569

$ws.Range("H136").Value = 1908.3793
$ws.Range("I136").Value = 741.05554
$ws.Range("J136").Value = 3818.5454
$ws.Range("K136").Value = 2223.16662
$ws.Range("L136").Value = 11455.6362
$ws.Range("M136").Value = 326.83338
$ws.Range("N136").Value = -16555.6362

$ws.Range("H140").Value = 35197.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 35197.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 35197.5
$ws.Range("N140").Value = -45557.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1675.7693
$ws.Range("I118").Value = 2368.1667
$ws.Range("J118").Value = 1082.2858
$ws.Range("K118").Value = 7104.500100000001
$ws.Range("L118").Value = 3246.8574
$ws.Range("M118").Value = -5861.500100000001
$ws.Range("N118").Value = -5732.857400000001

$ws.Range("H134").Value = 3172.3076
$ws.Range("I134").Value = 1809.4736
$ws.Range("J134").Value = 6871.4287
$ws.Range("K134").Value = 5428.4208
$ws.Range("L134").Value = 20614.2861
$ws.Range("M134").Value = -358.4207999999999
$ws.Range("N134").Value = -30754.2861

$ws.Range("H137").Value = 23817244
$ws.Range("I137").Value = 166666670
$ws.Range("J137").Value = 9005.5
$ws.Range("K137").Value = 500000010
$ws.Range("L137").Value = 27016.5
$ws.Range("M137").Value = -499994910
$ws.Range("N137").Value = -37216.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2668303.2
$ws.Range("I122").Value = 3176094
$ws.Range("J122").Value = 2402
$ws.Range("K122").Value = 9528282
$ws.Range("L122").Value = 7206
$ws.Range("M122").Value = -9525832
$ws.Range("N122").Value = -12106

$ws.Range("H132").Value = 6598.364
$ws.Range("I132").Value = 6701.4287
$ws.Range("J132").Value = 6550.2666
$ws.Range("K132").Value = 20104.2861
$ws.Range("L132").Value = 19650.7998
$ws.Range("M132").Value = -17574.2861
$ws.Range("N132").Value = -24710.7998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7950
$ws.Range("I7").Value = 19333.334
$ws.Range("J7").Value = 5323.077
$ws.Range("K7").Value = 19333.334
$ws.Range("L7").Value = 5323.077
$ws.Range("M7").Value = -19221.334
$ws.Range("N7").Value = -5547.077

$ws.Range("H122").Value = 4648.115
$ws.Range("I122").Value = 5203.7856
$ws.Range("J122").Value = 3999.8333
$ws.Range("K122").Value = 15611.3568
$ws.Range("L122").Value = 11999.4999
$ws.Range("M122").Value = -13161.3568
$ws.Range("N122").Value = -16899.4999

$ws.Range("H126").Value = 7950
$ws.Range("I126").Value = 19333.334
$ws.Range("J126").Value = 5323.077
$ws.Range("K126").Value = 58000.00199999999
$ws.Range("L126").Value = 15969.231
$ws.Range("M126").Value = -55530.00199999999
$ws.Range("N126").Value = -20909.231

$ws.Range("H135").Value = 134549.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 134549.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 134549.5
$ws.Range("N135").Value = -144689.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30249
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30249
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30249
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -31497

$ws.Range("H66").Value = 30249
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30249
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90747
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -96987

$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -31498

$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -97488

$ws.Range("H107").Value = 1289.2727
$ws.Range("I107").Value = 1810.8572
$ws.Range("J107").Value = 376.5
$ws.Range("K107").Value = 5432.571599999999
$ws.Range("L107").Value = 1129.5
$ws.Range("M107").Value = -3512.571599999999
$ws.Range("N107").Value = -4969.5

$ws.Range("H122").Value = 3002.5
$ws.Range("I122").Value = 2988
$ws.Range("J122").Value = 3031.5
$ws.Range("K122").Value = 8964
$ws.Range("L122").Value = 9094.5
$ws.Range("M122").Value = -6514
$ws.Range("N122").Value = -13994.5

$ws.Range("H136").Value = 1161.5333
$ws.Range("I136").Value = 1132.2174
$ws.Range("J136").Value = 1257.8572
$ws.Range("K136").Value = 3396.6522
$ws.Range("L136").Value = 3773.5716
$ws.Range("M136").Value = -846.6522
$ws.Range("N136").Value = -8873.571599999999
